$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 16988.902
$ws.Range("I129").Value = 572.13336
$ws.Range("J129").Value = 22228.299
$ws.Range("K129").Value = 1716.40008
$ws.Range("L129").Value = 66684.897
$ws.Range("M129").Value = 3283.59992
$ws.Range("N129").Value = -76684.897
$ws.Range("H137").Value = 1609.48
$ws.Range("I137").Value = 1049.2354
$ws.Range("J137").Value = 2800
$ws.Range("K137").Value = 3147.7062
$ws.Range("L137").Value = 8400
$ws.Range("M137").Value = -597.7062000000001
$ws.Range("N137").Value = -13500
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20105.838
$ws.Range("I32").Value = 20302.793
$ws.Range("J32").Value = 17250
$ws.Range("K32").Value = 20302.793
$ws.Range("L32").Value = 17250
$ws.Range("M32").Value = -20015.793
$ws.Range("N32").Value = -17824
$ws.Range("H97").Value = 496.91306
$ws.Range("I97").Value = 446.14285
$ws.Range("J97").Value = 1030
$ws.Range("K97").Value = 446.14285
$ws.Range("L97").Value = 1030
$ws.Range("M97").Value = 49.85714999999999
$ws.Range("N97").Value = -2022
$ws.Range("H122").Value = 1183.8
$ws.Range("I122").Value = 984
$ws.Range("J122").Value = 1650
$ws.Range("K122").Value = 2952
$ws.Range("L122").Value = 4950
$ws.Range("M122").Value = -502
$ws.Range("N122").Value = -9850
$ws.Range("H132").Value = 2815.7297
$ws.Range("I132").Value = 2827.1562
$ws.Range("K132").Value = 8481.4686
$ws.Range("M132").Value = -5951.4686
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1781.2222
$ws.Range("J86").Value = 2374.75
$ws.Range("L86").Value = 2374.75
$ws.Range("N86").Value = -4620.75
$ws.Range("H89").Value = 1781.2222
$ws.Range("J89").Value = 2374.75
$ws.Range("L89").Value = 11873.75
$ws.Range("N89").Value = -23105.75
$ws.Range("H107").Value = 866.7646999999999
$ws.Range("I107").Value = 740.1
$ws.Range("J107").Value = 1047.7142
$ws.Range("K107").Value = 740.1
$ws.Range("L107").Value = 1047.7142
$ws.Range("M107").Value = 1179.9
$ws.Range("N107").Value = -4887.7142
$ws.Range("H134").Value = 25241.215
$ws.Range("I134").Value = 44568.61
$ws.Range("J134").Value = 1844.8948
$ws.Range("K134").Value = 133705.83
$ws.Range("L134").Value = 5534.6844
$ws.Range("M134").Value = -131170.83
$ws.Range("N134").Value = -10604.6844
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5211277.5
$ws.Range("I31").Value = 2662.476
$ws.Range("J31").Value = 15154997
$ws.Range("K31").Value = 2662.476
$ws.Range("L31").Value = 15154997
$ws.Range("M31").Value = -2367.476
$ws.Range("N31").Value = -15155587
$ws.Range("H34").Value = 5211277.5
$ws.Range("I34").Value = 2662.476
$ws.Range("J34").Value = 15154997
$ws.Range("K34").Value = 2662.476
$ws.Range("L34").Value = 15154997
$ws.Range("M34").Value = -2460.476
$ws.Range("N34").Value = -15155401
$ws.Range("H105").Value = 701.1
$ws.Range("I105").Value = 600.7143
$ws.Range("J105").Value = 935.3333
$ws.Range("K105").Value = 600.7143
$ws.Range("L105").Value = 935.3333
$ws.Range("M105").Value = 1146.2857
$ws.Range("N105").Value = -4429.3333
$ws.Range("H134").Value = 1374.8334
$ws.Range("I134").Value = 888.6667
$ws.Range("J134").Value = 2833.3333
$ws.Range("K134").Value = 2666.0001
$ws.Range("L134").Value = 8499.999899999999
$ws.Range("M134").Value = -131.0001000000002
$ws.Range("N134").Value = -13569.9999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2540.1929
$ws.Range("J131").Value = 882.5814
$ws.Range("L131").Value = 2647.7442
$ws.Range("N131").Value = -12727.7442
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1128.6666
$ws.Range("I97").Value = 1047.7778
$ws.Range("J97").Value = 1250
$ws.Range("K97").Value = 1047.7778
$ws.Range("L97").Value = 1250
$ws.Range("M97").Value = -551.7778000000001
$ws.Range("N97").Value = -2242
$ws.Range("H107").Value = 286.5
$ws.Range("I107").Value = 340
$ws.Range("J107").Value = 249.46153
$ws.Range("K107").Value = 340
$ws.Range("L107").Value = 249.46153
$ws.Range("M107").Value = 1580
$ws.Range("N107").Value = -4089.46153
$ws.Range("H126").Value = 2441.1765
$ws.Range("I126").Value = 2775
$ws.Range("J126").Value = 1964.2858
$ws.Range("K126").Value = 8325
$ws.Range("L126").Value = 5892.857400000001
$ws.Range("M126").Value = -5855
$ws.Range("N126").Value = -10832.8574
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 35716256
$ws.Range("I7").Value = 1613.3334
$ws.Range("J7").Value = 62502236
$ws.Range("K7").Value = 1613.3334
$ws.Range("L7").Value = 62502236
$ws.Range("M7").Value = -1501.3334
$ws.Range("N7").Value = -62502460
$ws.Range("H16").Value = 548.5714
$ws.Range("I16").Value = 483.33334
$ws.Range("J16").Value = 597.5
$ws.Range("K16").Value = 483.33334
$ws.Range("L16").Value = 597.5
$ws.Range("M16").Value = -313.33334
$ws.Range("N16").Value = -937.5
$ws.Range("H40").Value = 1240.1177
$ws.Range("I40").Value = 1176.3077
$ws.Range("J40").Value = 1447.5
$ws.Range("K40").Value = 1176.3077
$ws.Range("L40").Value = 1447.5
$ws.Range("M40").Value = -1040.3077
$ws.Range("N40").Value = -1719.5
$ws.Range("H126").Value = 35716256
$ws.Range("I126").Value = 1613.3334
$ws.Range("J126").Value = 62502236
$ws.Range("K126").Value = 4840.0002
$ws.Range("L126").Value = 187506708
$ws.Range("M126").Value = -2370.0002
$ws.Range("N126").Value = -187511648
$ws.Range("H133").Value = 21669.555
$ws.Range("J133").Value = 21669.555
$ws.Range("L133").Value = 21669.555
$ws.Range("N133").Value = -26729.555
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4613.1665
$ws.Range("I81").Value = 4544.75
$ws.Range("J81").Value = 4750
$ws.Range("K81").Value = 9089.5
$ws.Range("L81").Value = 9500
$ws.Range("M81").Value = -8028.5
$ws.Range("N81").Value = -11622
$ws.Range("H84").Value = 4613.1665
$ws.Range("I84").Value = 4544.75
$ws.Range("J84").Value = 4750
$ws.Range("K84").Value = 45447.5
$ws.Range("L84").Value = 47500
$ws.Range("M84").Value = -40143.5
$ws.Range("N84").Value = -58108
$ws.Range("H122").Value = 52534.9
$ws.Range("I122").Value = 69047.336
$ws.Range("J122").Value = 2997.6
$ws.Range("K122").Value = 207142.008
$ws.Range("L122").Value = 8992.799999999999
$ws.Range("M122").Value = -204692.008
$ws.Range("N122").Value = -13892.8
$ws.Range("H126").Value = 83343520
$ws.Range("I126").Value = 166685920
$ws.Range("J126").Value = 1125.3334
$ws.Range("K126").Value = 500057760
$ws.Range("L126").Value = 3376.0002
$ws.Range("M126").Value = -500055290
$ws.Range("N126").Value = -8316.0002
